# Auto-generated Excel COM-interop script to update the cryptos price table
# Applies the diff: updates Price (D) and Volume(1h) (E) columns for rows 2-51,
# plus a few Coin/Link (B/C) swaps where the underlying coin order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '42.660.68'
Set-TextCell $ws 'E2' '  -7.66%  '

# Row 3
Set-TextCell $ws 'D3' '2.527.63'
Set-TextCell $ws 'E3' '  -2.62%  '

# Row 4
Set-TextCell $ws 'D4' '1.00'
Set-TextCell $ws 'E4' '  +0.03%  '

# Row 5
Set-TextCell $ws 'D5' '294.95'
Set-TextCell $ws 'E5' '  -4.19%  '

# Row 6
Set-TextCell $ws 'D6' '91.18'
Set-TextCell $ws 'E6' '  -7.97%  '

# Row 7
Set-TextCell $ws 'E7' '  -4.45%  '

# Row 8
Set-TextCell $ws 'E8' '  +0.04%  '

# Row 9
Set-TextCell $ws 'D9' '0.542'
Set-TextCell $ws 'E9' '  -6.30%  '

# Row 10
Set-TextCell $ws 'D10' '35.35'
Set-TextCell $ws 'E10' '  -9.31%  '

# Row 11
Set-TextCell $ws 'D11' '0.0801'
Set-TextCell $ws 'E11' '  -4.52%  '

# Row 12
Set-TextCell $ws 'D12' '7.60'
Set-TextCell $ws 'E12' '  -6.55%  '

# Row 13
Set-TextCell $ws 'D13' '0.106'
Set-TextCell $ws 'E13' '  +0.31%  '

# Row 14
Set-TextCell $ws 'D14' '2.912.14'
Set-TextCell $ws 'E14' '  -2.75%  '

# Row 15
Set-TextCell $ws 'D15' '2.523.96'
Set-TextCell $ws 'E15' '  -2.79%  '

# Row 16
Set-TextCell $ws 'E16' '  -6.34%  '

# Row 17
Set-TextCell $ws 'D17' '14.03'
Set-TextCell $ws 'E17' '  -5.66%  '

# Row 18
Set-TextCell $ws 'D18' '42.730.18'
Set-TextCell $ws 'E18' '  -7.70%  '

# Row 19
Set-TextCell $ws 'B19' 'ShibaInu'
Set-TextCell $ws 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D19' '0.0₃0959'
Set-TextCell $ws 'E19' '  -5.03%  '

# Row 20
Set-TextCell $ws 'B20' 'Uniswap'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws 'D20' '6.50'
Set-TextCell $ws 'E20' '  -3.10%  '

# Row 21
Set-TextCell $ws 'D21' '12.23'
Set-TextCell $ws 'E21' '  -4.57%  '

# Row 22
Set-TextCell $ws 'D22' '71.37'
Set-TextCell $ws 'E22' '  -0.20%  '

# Row 23
Set-TextCell $ws 'D23' '256.65'
Set-TextCell $ws 'E23' '  -6.32%  '

# Row 24
Set-TextCell $ws 'D24' '2.87'
Set-TextCell $ws 'E24' '  -5.11%  '

# Row 25
Set-TextCell $ws 'E25' '  -4.18%  '

# Row 26
Set-TextCell $ws 'D26' '28.59'
Set-TextCell $ws 'E26' '  -3.45%  '

# Row 27
Set-TextCell $ws 'E27' '  +0.14%  '

# Row 28
Set-TextCell $ws 'B28' 'Cosmos'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D28' '9.86'
Set-TextCell $ws 'E28' '  -8.44%  '

# Row 29
Set-TextCell $ws 'B29' 'Toncoin'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D29' '2.14'
Set-TextCell $ws 'E29' '  -3.24%  '

# Row 30
Set-TextCell $ws 'D30' '35.87'
Set-TextCell $ws 'E30' '  -6.05%  '

# Row 31
Set-TextCell $ws 'D31' '5.87'
Set-TextCell $ws 'E31' '  -6.28%  '

# Row 32
Set-TextCell $ws 'B32' 'Monero'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D32' '150.20'
Set-TextCell $ws 'E32' '  -3.24%  '

# Row 33
Set-TextCell $ws 'B33' 'LidoDAOToken'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D33' '3.39'
Set-TextCell $ws 'E33' '  -4.84%  '

# Row 34
Set-TextCell $ws 'D34' '2.14'
Set-TextCell $ws 'E34' '  -4.10%  '

# Row 35
Set-TextCell $ws 'D35' '2.69'
Set-TextCell $ws 'E35' '  -3.31%  '

# Row 36
Set-TextCell $ws 'D36' '0.0788'
Set-TextCell $ws 'E36' '  -5.42%  '

# Row 37
Set-TextCell $ws 'D37' '0.113'
Set-TextCell $ws 'E37' '  -7.16%  '

# Row 38
Set-TextCell $ws 'D38' '0.119'
Set-TextCell $ws 'E38' '  -3.69%  '

# Row 39
Set-TextCell $ws 'D39' '23.52'
Set-TextCell $ws 'E39' '  +5.13%  '

# Row 40
Set-TextCell $ws 'D40' '16.08'
Set-TextCell $ws 'E40' '  +1.21%  '

# Row 41
Set-TextCell $ws 'D41' '3.37'
Set-TextCell $ws 'E41' '  -5.50%  '

# Row 42
Set-TextCell $ws 'E42' '  -7.19%  '

# Row 43
Set-TextCell $ws 'D43' '3.77'
Set-TextCell $ws 'E43' '  -4.98%  '

# Row 44
Set-TextCell $ws 'D44' '2.046.11'
Set-TextCell $ws 'E44' '  -3.35%  '

# Row 45
Set-TextCell $ws 'E45' '  +0.03%  '

# Row 46
Set-TextCell $ws 'D46' '83.84'
Set-TextCell $ws 'E46' '  -11.85%  '

# Row 47
Set-TextCell $ws 'D47' '1.59'
Set-TextCell $ws 'E47' '  +2.92%  '

# Row 48
Set-TextCell $ws 'D48' '8.77'
Set-TextCell $ws 'E48' '  -8.34%  '

# Row 49
Set-TextCell $ws 'D49' '2.768.81'
Set-TextCell $ws 'E49' '  -2.82%  '

# Row 50
Set-TextCell $ws 'E50' '  -4.99%  '

# Row 51
Set-TextCell $ws 'D51' '102.62'
Set-TextCell $ws 'E51' '  -5.57%  '
